# Update version to v0.4
# - Insert a new column before D (shifts old D -> E, new D inherits format from column C)
# - Change C2 from the JSON string to the plain id locator string
# - Set the new D3 cell to the JSON string (plain text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; existing column D (and its contents/format) moves to E.
$ws.Range("D1").EntireColumn.Insert()

# Row 1 header: new column D mirrors the "mouseOver" header (same as C1)
$ws.Range("D1").Value = "mouseOver"

# Row 2: C2 becomes the plain locator string instead of the JSON blob
$ws.Range("C2").Value = "id=btn1"

# Row 3: new D3 holds the JSON payload that used to live in C2
$ws.Range("D3").Value = '{"target": "id=btn1"}'

# Match the author's final selection / column width state
$ws.Range("E5").Select() | Out-Null
$ws.Columns("D").ColumnWidth = 18.95
